# Refresh the cryptocurrency price/volume snapshot (columns D and E).
# The source values are plain text, not locale-aware numbers - e.g.
# prices such as "64.924.10" use '.' as a thousands separator and
# would otherwise be mis-parsed - so any new value that looks like a
# plain number is written with a leading apostrophe. That forces Excel
# to store it as text, matching the original inline-string cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.924.10"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "3.394.86"
$ws.Range("E3").Value = "  +1.23%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'560.38"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").Value = "'174.95"
$ws.Range("E6").Value = "  +1.55%  "
$ws.Range("D7").Value = "'0.627"
$ws.Range("E7").Value = "  +2.31%  "
$ws.Range("D8").Value = "3.386.33"
$ws.Range("E8").Value = "  +1.36%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  +10.61%  "
$ws.Range("D11").Value = "'0.633"
$ws.Range("E11").Value = "  +2.05%  "
$ws.Range("D12").Value = "'54.76"
$ws.Range("E12").Value = "  +1.44%  "
$ws.Range("E13").Value = "  +4.50%  "
$ws.Range("D14").Value = "'9.15"
$ws.Range("E14").Value = "  +2.60%  "
$ws.Range("D15").Value = "3.933.39"
$ws.Range("E15").Value = "  +1.51%  "
$ws.Range("E16").Value = "  +2.73%  "
$ws.Range("D17").Value = "3.394.88"
$ws.Range("E17").Value = "  +1.15%  "
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").Value = "'11.91"
$ws.Range("E19").Value = "  +1.44%  "
$ws.Range("D20").Value = "64.870.94"
$ws.Range("E20").Value = "  +1.42%  "
$ws.Range("E21").Value = "  +2.17%  "
$ws.Range("D22").Value = "'468.17"
$ws.Range("E22").Value = "  +15.97%  "
$ws.Range("D23").Value = "'4.99"
$ws.Range("E23").Value = "  +16.73%  "
$ws.Range("D24").Value = "'4.14"
$ws.Range("E24").Value = "  +1.28%  "
$ws.Range("D25").Value = "'86.33"
$ws.Range("E25").Value = "  +4.15%  "
$ws.Range("D26").Value = "'13.61"
$ws.Range("E26").Value = "  +2.21%  "
$ws.Range("D27").Value = "'10.86"
$ws.Range("E27").Value = "  +1.55%  "
$ws.Range("E28").Value = "  +4.03%  "
$ws.Range("E29").Value = "  +1.46%  "
$ws.Range("D30").Value = "'30.63"
$ws.Range("E30").Value = "  +4.75%  "
$ws.Range("E31").Value = "  +4.32%  "
$ws.Range("D32").Value = "'11.53"
$ws.Range("E32").Value = "  +1.39%  "
$ws.Range("D33").Value = "'584.16"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  +2.14%  "
$ws.Range("D35").Value = "'60.07"
$ws.Range("E35").Value = "  +3.57%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("E37").Value = "  -3.84%  "
$ws.Range("D38").Value = "'35.91"
$ws.Range("E38").Value = "  +0.47%  "
$ws.Range("D39").Value = "'3.49"
$ws.Range("E39").Value = "  +0.93%  "
$ws.Range("D40").Value = "0.0₃0756"
$ws.Range("E40").Value = "  +1.75%  "
$ws.Range("E41").Value = "  +0.94%  "
$ws.Range("D42").Value = "3.096.51"
$ws.Range("E42").Value = "  -1.67%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("D44").Value = "'2.86"
$ws.Range("E44").Value = "  +0.37%  "
$ws.Range("E45").Value = "  +1.99%  "
$ws.Range("E46").Value = "  +2.40%  "
$ws.Range("E47").Value = "  -1.58%  "
$ws.Range("E48").Value = "  +4.20%  "
$ws.Range("E49").Value = "  -2.91%  "
$ws.Range("D50").Value = "'138.12"
$ws.Range("E50").Value = "  +4.08%  "
$ws.Range("D51").Value = "'8.40"
$ws.Range("E51").Value = "  +3.58%  "
